$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 572 (shifts existing rows 572:647 down to 575:650).
# Excel's Insert copies formatting from the row above, which preserves the
# date number-format already present on column D.
$ws.Rows("572:574").Insert()

# New row 572 (week 2023-09-11, Agricola del Norte S.A. de Arica - Zanahoria, Primera)
$ws.Cells.Item(572, 1).Value = 1
$ws.Cells.Item(572, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(572, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(572, 4).Value = 45180
$ws.Cells.Item(572, 5).Value = 15
$ws.Cells.Item(572, 6).Value = 100114013
$ws.Cells.Item(572, 7).Value = "Zanahoria"
$ws.Cells.Item(572, 8).Value = "Sin especificar"
$ws.Cells.Item(572, 9).Value = "Primera"
$ws.Cells.Item(572, 10).Value = 35
$ws.Cells.Item(572, 11).Value = 7000
$ws.Cells.Item(572, 12).Value = 8000
$ws.Cells.Item(572, 13).Value = 7429
$ws.Cells.Item(572, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(572, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(572, 16).Value = 297
$ws.Cells.Item(572, 17).Value = 25
$ws.Cells.Item(572, 18).Value = "Hortaliza"

# New row 573 (same week, Segunda)
$ws.Cells.Item(573, 1).Value = 1
$ws.Cells.Item(573, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(573, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(573, 4).Value = 45180
$ws.Cells.Item(573, 5).Value = 15
$ws.Cells.Item(573, 6).Value = 100114013
$ws.Cells.Item(573, 7).Value = "Zanahoria"
$ws.Cells.Item(573, 8).Value = "Sin especificar"
$ws.Cells.Item(573, 9).Value = "Segunda"
$ws.Cells.Item(573, 10).Value = 25
$ws.Cells.Item(573, 11).Value = 5000
$ws.Cells.Item(573, 12).Value = 6000
$ws.Cells.Item(573, 13).Value = 5600
$ws.Cells.Item(573, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(573, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(573, 16).Value = 224
$ws.Cells.Item(573, 17).Value = 25
$ws.Cells.Item(573, 18).Value = "Hortaliza"

# New row 574 (same week, Tercera)
$ws.Cells.Item(574, 1).Value = 1
$ws.Cells.Item(574, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(574, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(574, 4).Value = 45180
$ws.Cells.Item(574, 5).Value = 15
$ws.Cells.Item(574, 6).Value = 100114013
$ws.Cells.Item(574, 7).Value = "Zanahoria"
$ws.Cells.Item(574, 8).Value = "Sin especificar"
$ws.Cells.Item(574, 9).Value = "Tercera"
$ws.Cells.Item(574, 10).Value = 20
$ws.Cells.Item(574, 11).Value = 3000
$ws.Cells.Item(574, 12).Value = 4000
$ws.Cells.Item(574, 13).Value = 3500
$ws.Cells.Item(574, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(574, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(574, 16).Value = 140
$ws.Cells.Item(574, 17).Value = 25
$ws.Cells.Item(574, 18).Value = "Hortaliza"
